$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 367, pushing existing rows
# (367-418) down to (369-420), preserving their formatting/values as-is.
$ws.Rows("367:368").Insert()

# New row 367 ("Primera")
$ws.Range("A367").Value = 7
$ws.Range("B367").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C367").Value = "Ñuble"
$ws.Range("D367").Value = "2023-08-04"
$ws.Range("E367").Value = 16
$ws.Range("F367").Value = 100112017
$ws.Range("G367").Value = "Apio"
$ws.Range("H367").Value = "Americana (o)"
$ws.Range("I367").Value = "Primera"
$ws.Range("J367").Value = 120
$ws.Range("K367").Value = 6000
$ws.Range("L367").Value = 6000
$ws.Range("M367").Value = 6000
$ws.Range("N367").Value = "$/docena de matas"
$ws.Range("O367").Value = "Provincia del Elquí"
$ws.Range("P367").Value = 1000
$ws.Range("Q367").Value = 6
$ws.Range("R367").Value = "Hortaliza"

# New row 368 ("Segunda")
$ws.Range("A368").Value = 7
$ws.Range("B368").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C368").Value = "Ñuble"
$ws.Range("D368").Value = "2023-08-04"
$ws.Range("E368").Value = 16
$ws.Range("F368").Value = 100112017
$ws.Range("G368").Value = "Apio"
$ws.Range("H368").Value = "Americana (o)"
$ws.Range("I368").Value = "Segunda"
$ws.Range("J368").Value = 80
$ws.Range("K368").Value = 5000
$ws.Range("L368").Value = 5000
$ws.Range("M368").Value = 5000
$ws.Range("N368").Value = "$/docena de matas"
$ws.Range("O368").Value = "Provincia del Elquí"
$ws.Range("P368").Value = 833
$ws.Range("Q368").Value = 6
$ws.Range("R368").Value = "Hortaliza"
